# Comprehension scores.xlsx - "Add files via upload"
# A new week (week 21) of data was appended to Sheet1 as row 22:
#   A22 = 21                         (week number)
#   B22 = 1.8798842592592593         (time spent that week, [h]:mm:ss like the rows above)
#   C22 = =SUM(B2:B22)+1.2708333333  (running total formula, same pattern as C2:C21)
#   D22 = the new "what I did" note for that week
# The selection/scroll position also moved to reflect the newly active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (week 21)
$ws.Range("A22").Value = 21

$ws.Range("B22").Value = 1.8798842592592593
$ws.Range("B22").NumberFormat = $ws.Range("B20").NumberFormat

$ws.Range("C22").Formula = "=SUM(B2:B22)+1.2708333333"
$ws.Range("C22").NumberFormat = $ws.Range("C21").NumberFormat

$ws.Range("D22").Value = "Oscuro deseo (Audiovisual, Spanish, Re-watch):32; [¿Cuáles son las MEJORES papas picantes de la tiendita? ](https://youtu.be/VbJMWqfeX1E) (Audiovisual, Spanish, New):41; "

# Reflect the scrolled/selected state from the saved file (best effort).
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("C23").Select()
